# Insert a new weekly record for "Feria Lagunitas de Puerto Montt - Haba" at row 54,
# pushing the existing rows 54-90 down to 55-91 (data itself is unchanged, only its
# row position shifts down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 54; everything currently at/after row 54 moves down one row.
$ws.Rows("54:54").Insert()

# Populate the newly inserted row 54 with the new weekly price record.
$ws.Range("A54").Value = 4
$ws.Range("B54").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C54").Value = "Los Lagos"
$ws.Range("D54").Value = 44762
$ws.Range("E54").Value = 10
$ws.Range("F54").Value = 100112026
$ws.Range("G54").Value = "Haba"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 80
$ws.Range("K54").Value = 21000
$ws.Range("L54").Value = 21000
$ws.Range("M54").Value = 21000
$ws.Range("N54").Value = "$/saco 25 kilos"
$ws.Range("O54").Value = "Provincia de Limarí"
$ws.Range("P54").Value = 840
$ws.Range("Q54").Value = 25
$ws.Range("R54").Value = "Hortaliza"
